$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Serpinf1"
$ws.Range("C2").Value = "Plxdc2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.220712
$ws.Range("H2").Value = 6.662135999999999
$ws.Range("I2").Value = 0.004164179109543329
$ws.Range("J2").Value = 0.00416417910954333
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.086735333333333
$ws.Range("N2").Value = 3.260206
$ws.Range("O2").Value = 0.007520028150622985
$ws.Range("P2").Value = 0.007520028150622985
$ws.Range("Q2").Value = 2.413326195557333
$ws.Range("R2").Value = 21.719935760016
$ws.Range("S2").Value = 0.00003131474412800199
$ws.Range("T2").Value = 0.00003131474412800199

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Serpinf1"
$ws.Range("C3").Value = "Plxdc2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.220712
$ws.Range("H3").Value = 6.662135999999999
$ws.Range("I3").Value = 0.004164179109543329
$ws.Range("J3").Value = 0.00416417910954333
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 124.5345636666667
$ws.Range("N3").Value = 373.603691
$ws.Range("O3").Value = 0.8617585126512408
$ws.Range("P3").Value = 0.8617585126512408
$ws.Range("Q3").Value = 276.5553999493306
$ws.Range("R3").Value = 2488.998599543976
$ws.Range("S3").Value = 0.003588516795853427
$ws.Range("T3").Value = 0.003588516795853428

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Serpinf1"
$ws.Range("C4").Value = "Plxdc2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.220712
$ws.Range("H4").Value = 6.662135999999999
$ws.Range("I4").Value = 0.004164179109543329
$ws.Range("J4").Value = 0.00416417910954333
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 18.89083733333333
$ws.Range("N4").Value = 56.672512
$ws.Range("O4").Value = 0.1307214591981362
$ws.Range("P4").Value = 0.1307214591981362
$ws.Range("Q4").Value = 41.95110915618133
$ws.Range("R4").Value = 377.559982405632
$ws.Range("S4").Value = 0.0005443475695618995
$ws.Range("T4").Value = 0.0005443475695618995

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Serpinf1"
$ws.Range("C5").Value = "Plxdc2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 527.8012189999999
$ws.Range("H5").Value = 1583.403657
$ws.Range("I5").Value = 0.9897090708526379
$ws.Range("J5").Value = 0.9897090708526382
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.086735333333333
$ws.Range("N5").Value = 3.260206
$ws.Range("O5").Value = 0.007520028150622985
$ws.Range("P5").Value = 0.007520028150622985
$ws.Range("Q5").Value = 573.5802336637047
$ws.Range("R5").Value = 5162.222102973342
$ws.Range("S5").Value = 0.007442640073738756
$ws.Range("T5").Value = 0.007442640073738758

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Serpinf1"
$ws.Range("C6").Value = "Plxdc2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 527.8012189999999
$ws.Range("H6").Value = 1583.403657
$ws.Range("I6").Value = 0.9897090708526379
$ws.Range("J6").Value = 0.9897090708526382
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 124.5345636666667
$ws.Range("N6").Value = 373.603691
$ws.Range("O6").Value = 0.8617585126512408
$ws.Range("P6").Value = 0.8617585126512408
$ws.Range("Q6").Value = 65729.49451089978
$ws.Range("R6").Value = 591565.450598098
$ws.Range("S6").Value = 0.8528902168554108
$ws.Range("T6").Value = 0.852890216855411

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Serpinf1"
$ws.Range("C7").Value = "Plxdc2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 527.8012189999999
$ws.Range("H7").Value = 1583.403657
$ws.Range("I7").Value = 0.9897090708526379
$ws.Range("J7").Value = 0.9897090708526382
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 18.89083733333333
$ws.Range("N7").Value = 56.672512
$ws.Range("O7").Value = 0.1307214591981362
$ws.Range("P7").Value = 0.1307214591981362
$ws.Range("Q7").Value = 9970.606972464042
$ws.Range("R7").Value = 89735.46275217638
$ws.Range("S7").Value = 0.1293762139234884
$ws.Range("T7").Value = 0.1293762139234884

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Serpinf1"
$ws.Range("C8").Value = "Plxdc2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.26733
$ws.Range("H8").Value = 9.80199
$ws.Range("I8").Value = 0.006126750037818593
$ws.Range("J8").Value = 0.006126750037818595
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.086735333333333
$ws.Range("N8").Value = 3.260206
$ws.Range("O8").Value = 0.007520028150622985
$ws.Range("P8").Value = 0.007520028150622985
$ws.Range("Q8").Value = 3.55072295666
$ws.Range("R8").Value = 31.95650660994
$ws.Range("S8").Value = 0.00004607333275622626
$ws.Range("T8").Value = 0.00004607333275622627

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Serpinf1"
$ws.Range("C9").Value = "Plxdc2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.26733
$ws.Range("H9").Value = 9.80199
$ws.Range("I9").Value = 0.006126750037818593
$ws.Range("J9").Value = 0.006126750037818595
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 124.5345636666667
$ws.Range("N9").Value = 373.603691
$ws.Range("O9").Value = 0.8617585126512408
$ws.Range("P9").Value = 0.8617585126512408
$ws.Range("Q9").Value = 406.89551590501
$ws.Range("R9").Value = 3662.05964314509
$ws.Range("S9").Value = 0.005279778999976484
$ws.Range("T9").Value = 0.005279778999976486

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Serpinf1"
$ws.Range("C10").Value = "Plxdc2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.26733
$ws.Range("H10").Value = 9.80199
$ws.Range("I10").Value = 0.006126750037818593
$ws.Range("J10").Value = 0.006126750037818595
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 18.89083733333333
$ws.Range("N10").Value = 56.672512
$ws.Range("O10").Value = 0.1307214591981362
$ws.Range("P10").Value = 0.1307214591981362
$ws.Range("Q10").Value = 61.72259954432
$ws.Range("R10").Value = 555.5033958988799
$ws.Range("S10").Value = 0.0008008977050858829
$ws.Range("T10").Value = 0.000800897705085883

Write-Host "Edit complete"
